$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "absence"
$ws.Range("B2").Value = 30

$ws.Range("A3").Value = "présence"
$ws.Range("B3").Value = 63
